$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new journal entry text (matches the shared-string added by the diff)
$newText = @"
J'ai modifié la base de données en ajoutant un champs pour les illustrations des chaussures, 
j'ai ajouté des données pour que dans la page de shopping il y ait des chaussures qui s'affichent. J'ai réglé un problème que j'avais lorsque je m'inscrivais sur le site, il y avait une erreur PhP qui disait qu'il ne connaissait pas des variables.
"@
$newText = $newText.TrimEnd("`r", "`n")

# Copy the formatting (number format / wrap / fill etc.) from the last
# existing data row (21) down onto the new row (22), mirroring the
# look of the prior journal entries.
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)

$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)

$ws.Range("C21").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Fill in the new journal entry (date 2018-02-14 / serial 43145)
$ws.Range("A22").Value = 43145
$ws.Range("B22").Value = $newText
$ws.Range("C22").Value = "3 périodes"

# Row 22 needs to be tall enough for the wrapped, multi-line text.
$ws.Rows.Item(22).RowHeight = 75

# Column C got an explicit best-fit width once data was added to it
# beyond the header row (target stored width ~10.29 chars).
$ws.Columns.Item(3).ColumnWidth = 9.5

# Move the active selection to C23 (one row below the newly added row),
# matching where the cursor ends up after entering the new row's data.
$ws.Range("C23").Select()
